$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 updates
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 1.2848039557918909
$ws.Range("D2").Value = 0.30370967635272378
$ws.Range("E2").Value = 3.854618248298646

# Row 3 updates
$ws.Range("B3").Value = 0.79310944099776215
$ws.Range("C3").Value = 2.0886367795514373
$ws.Range("D3").Value = 2.0790317919247618
$ws.Range("E3").Value = 3.3165995798502315

# Update selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
